$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-10-18 Saturday" "2025-10-19 Sunday"

Replace-Text "45÷8=5, 5" "94÷4=23, 2"
Replace-Text "14÷8=1, 6" "84÷4=21, 0"
Replace-Text "19÷2=9, 1" "90÷5=18, 0"
Replace-Text "59÷6=9, 5" "78÷5=15, 3"
Replace-Text "75÷7=10, 5" "67÷5=13, 2"

Replace-Text "79÷2=39, 1" "33÷8=4, 1"
Replace-Text "76÷6=12, 4" "10÷7=1, 3"
Replace-Text "22÷6=3, 4" "49÷7=7, 0"
Replace-Text "99÷9=11, 0" "92÷8=11, 4"
Replace-Text "65÷8=8, 1" "56÷7=8, 0"

Replace-Text "47÷8=5, 7" "88÷5=17, 3"
Replace-Text "18÷4=4, 2" "80÷4=20, 0"
Replace-Text "74÷6=12, 2" "75÷4=18, 3"
Replace-Text "21÷6=3, 3" "42÷6=7, 0"
Replace-Text "26÷9=2, 8" "65÷4=16, 1"

Replace-Text "62÷5=12, 2" "36÷7=5, 1"
Replace-Text "42÷3=14, 0" "40÷7=5, 5"
Replace-Text "14÷9=1, 5" "11÷2=5, 1"
Replace-Text "17÷4=4, 1" "41÷5=8, 1"
Replace-Text "70÷9=7, 7" "38÷9=4, 2"

Replace-Text "18÷6=3, 0" "71÷8=8, 7"
Replace-Text "93÷4=23, 1" "13÷8=1, 5"
Replace-Text "60÷3=20, 0" "74÷5=14, 4"
Replace-Text "65÷9=7, 2" "15÷7=2, 1"
Replace-Text "82÷4=20, 2" "80÷7=11, 3"
